# paises.xlsx -- "Update countries & provincias Spain"
#
# Refreshes the daily COVID-19 country table on sheet "Pais":
#   - bumps the "last updated" banner in A1 to 18:57
#   - writes the new case/death/recovered counters pulled in for this
#     update (several countries' row data changed)
#   - two country-name/rank swaps fell out of the refreshed totals:
#       * Irak overtook Turquia (rows 21/22 swap identity + data)
#       * Birmania was re-sorted to just after Sri Lanka, pushing
#         Trinidad yTobago / Guadalupe / Aruba down one row each
#         (rows 136-139)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 18:57"

# Each entry: row, country label (col A), then B,C,D,E,F,G,H
$rows = @(
    @(4,   "Estados Unidos",    6721465, 13007, 3982775, 2540010, 0, 160, 198680),
    @(6,   "Brasil",            4335066,  4611, 3573958,  629372, 0,  73, 131736),
    @(21,  "Irak",               294478,  4169,  229132,   57260, 0,  72,   8086),
    @(22,  "Turquia",            292878,  1716,  260058,   25701, 0,  63,   7119),
    @(25,  "Alemania",           262688,  1390,  235700,   17557, 0,   3,   9431),
    @(32,  "Ecuador",            118911,   317,   97063,   10926, 0,  19,  10922),
    @(58,  "Nepal",               55329,  1170,   39576,   15393, 0,  15,    360),
    @(59,  "Argelia",             48496,   242,   34204,   12672, 0,   8,   1620),
    @(68,  "Chequia",             36722,   534,   22020,   14237, 0,   9,    465),
    @(77,  "Libano",              25401,  1091,    8765,   16390, 0,   5,    246),
    @(91,  "Grecia",              13420,   180,    3804,    9306, 0,   5,    310),
    @(99,  "Guayana Francesa",     9552,    31,    9156,     333, 0,   0,     63),
    @(136, "Birmania",             3195,   263,     790,    2373, 0,  12,     32),
    @(137, "Trinidad yTobago",     3091,    49,     787,    2250, 0,   1,     54),
    @(138, "Guadalupe",            3080,     0,     837,    2219, 0,   0,     24),
    @(139, "Aruba",                3046,     0,    1542,    1486, 0,   0,     18),
    @(141, "Mali",                 2935,    11,    2289,     518, 0,   0,    128),
    @(182, "Isla de Man",           339,     2,     312,       3, 0,   0,     24)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Range("A$row").Value = $r[1]
    $ws.Range("B$row").Value = $r[2]
    $ws.Range("C$row").Value = $r[3]
    $ws.Range("D$row").Value = $r[4]
    $ws.Range("E$row").Value = $r[5]
    $ws.Range("F$row").Value = $r[6]
    $ws.Range("G$row").Value = $r[7]
    $ws.Range("H$row").Value = $r[8]
}
